# Scene.xlsx fix: "fixed for some error configure file"
#
# 1) M6 (CanClone for the "SelectScene/City" row) was mistakenly left at 0;
#    correct it to 1.
# 2) Leave the sheet's cursor/selection parked on N12 (post-edit cleanup
#    position), matching where the editor ended up after making the fix.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M6").Value = 1

$ws.Range("N12").Select()
